# Regenerate save_data column G ("K") values: use K (strikeouts) instead of
# the old Strike# metric, after recalculating std/mean and writing s_vals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new value for column G ("K")
$kValues = @{
    2  = 2
    3  = 5
    4  = 6
    5  = 0
    6  = 1
    7  = 9
    8  = 5
    9  = 4
    10 = 6
    11 = 4
    12 = 1
    13 = 3
    14 = 2
    15 = 2
    16 = 3
    17 = 7
    18 = 4
    19 = 7
    20 = 2
    21 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
